# Add data for 2022-07-26: the "through" date of the running July 2022
# month-to-date column moves from July 17 to July 18, and each neighborhood
# that logged a carjacking on "day 18" of July gets +1 in the corresponding
# "July <year>" column (column B = July 2022, and the other July columns for
# each prior year: I=2021, P=2020, W=2019, AD=2018, AK=2017, AR=2016).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the month-to-date label/header cell.
$ws.Name = "Through 2022-07-18"
$ws.Range("B1").Value = "July 2022 (through July 18)"

# Austin
$ws.Range("I2").Value = 9
$ws.Range("AD2").Value = 8

# Englewood
$ws.Range("I3").Value = 4

# Garfield Park
$ws.Range("B5").Value = 6
$ws.Range("AD5").Value = 4

# Grand Crossing
$ws.Range("P6").Value = 3

# North Lawndale
$ws.Range("I8").Value = 6

# Douglas
$ws.Range("I15").Value = 3

# Washington Heights
$ws.Range("P16").Value = 1
$ws.Range("AD16").Value = 2

# West Loop
$ws.Range("I18").Value = 3

# South Shore
$ws.Range("I19").Value = 4

# South Deering
$ws.Range("AR24").Value = 1

# Fuller Park
$ws.Range("P31").Value = 1

# West Town
$ws.Range("B38").Value = 3
$ws.Range("I38").Value = 5

# Logan Square
$ws.Range("AR41").Value = 2

# New City
$ws.Range("AK44").Value = 1
$ws.Range("AR44").Value = 2

# Little Italy, UIC
$ws.Range("AD47").Value = 2

# Grand Boulevard
$ws.Range("B49").Value = 1

# Chatham
$ws.Range("P52").Value = 9

# Albany Park
$ws.Range("I58").Value = 2
$ws.Range("AR58").Value = 3

# Hermosa
$ws.Range("AR76").Value = 1

# Lake View
$ws.Range("W78").Value = 1
